$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 49
$ws1.Range("F4").Value = 835
$ws1.Range("F7").Value = 592
$ws1.Range("F9").Value = 57
$ws1.Range("F12").Value = 659
$ws1.Range("F14").Value = 1809
$ws1.Range("F15").Value = 352
$ws1.Range("F16").Value = 3153
$ws1.Range("F17").Value = 320
$ws1.Range("F19").Value = 53
$ws1.Range("F21").Value = 130

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5329
$ws3.Range("F3").Value = 322
$ws3.Range("F4").Value = 258

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 49
$ws4.Range("F3").Value = 5329
$ws4.Range("F4").Value = 322
$ws4.Range("F6").Value = 258
$ws4.Range("F14").Value = 835
$ws4.Range("F19").Value = 592
$ws4.Range("F21").Value = 57
$ws4.Range("F27").Value = 659
$ws4.Range("F30").Value = 1809
$ws4.Range("F31").Value = 352
$ws4.Range("F32").Value = 3153
$ws4.Range("F34").Value = 320
$ws4.Range("F36").Value = 53
$ws4.Range("F40").Value = 130
